# Use BAU BECCS ban rather than policy lever to remove it from wedge diagram.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BBNPPTY")

# Row 21 corresponds to "biomass w CCS" (BECCS). Set the annual ban flags
# (columns B:AE, years 2021-2050) from 0 to 1.
$ws.Range("B21:AE21").Value = 1

# Reflect the selection left after making this edit in the UI.
$ws.Activate()
$ws.Range("B21:AE21").Select()
